$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.77294921875
$ws.Range("B1").Value = 6.162662982940674
$ws.Range("C1").Value = 8.474946975708008
$ws.Range("D1").Value = 7.040911674499512
$ws.Range("E1").Value = 2.537656784057617
